# "change in the pff history..." -- re-sync the Backlog Task workbook:
#  - re-merge the per-row header cells on "Feature Source" so the stored
#    merge order matches a later resave (rows 1-6, then 7-18, then 19-27,
#    each block re-merged "last row first, then ascending")
#  - scroll/select "Test Summary" down to the bottom block and update the
#    three running totals in column C (rows 35-37)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Feature Source" - re-merge A1:D1 .. A27:D27 in three row-blocks so
#    the mergeCells list is rewritten in the same grouped order as the
#    target file. Unmerging/merging a row resets its cell borders in
#    this engine, so every touched row's formatting is restored right
#    afterwards from a pristine copy stashed outside the used range.
# ---------------------------------------------------------------------
$wsFeature = $wb.Worksheets.Item("Feature Source")

# Stash a pristine copy of the uniform row formatting far outside the
# table (so later fix-ups never depend on a row we are about to touch).
$wsFeature.Range("A2:D2").Copy()
$wsFeature.Range("A100:D100").PasteSpecial(-4122)
$excel.CutCopyMode = $false

function Remerge-Row($sheet, $row) {
    $rng = $sheet.Range("A" + $row + ":D" + $row)
    $rng.UnMerge()
    $rng.Merge()
    $sheet.Range("A100:D100").Copy()
    $rng.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# Block rows 1-6: A6 re-merged first, then A1..A5 ascending
Remerge-Row $wsFeature 6
for ($r = 1; $r -le 5; $r++) { Remerge-Row $wsFeature $r }

# Block rows 7-18: A18 re-merged first, then A7..A17 ascending
Remerge-Row $wsFeature 18
for ($r = 7; $r -le 17; $r++) { Remerge-Row $wsFeature $r }

# Block rows 19-27: A25,A26,A27 re-merged first, then A19..A24 ascending
Remerge-Row $wsFeature 25
Remerge-Row $wsFeature 26
Remerge-Row $wsFeature 27
for ($r = 19; $r -le 24; $r++) { Remerge-Row $wsFeature $r }

# Remove the scratch formatting helper row
$wsFeature.Range("A100:D100").Clear()

# ---------------------------------------------------------------------
# 2) "Test Summary" - update the running totals and move the view/
#    selection down to the last block of the sheet.
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Test Summary")

$wsSummary.Range("C35").Value = 28
$wsSummary.Range("C36").Value = 34
$wsSummary.Range("C37").Value = 20

$wsSummary.Activate()
$excel.ActiveWindow.ScrollRow = 25
[void]$wsSummary.Range("J36").Select()
